$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "ERF(X,Y)"
$ws.Range("C25").Formula = "=ERF(A25,B25)"
$ws.Range("C26").Formula = "=ERF(A26,B26)"

$ws.Range("C24").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
